$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2022-05-30 T 21:36:43 UTC"
$ws.Range("B2").Value = 30803.0368329025
$ws.Range("C2").Value = 1.013199
$ws.Range("D2").Value = 1.204574

$ws.Range("A3").Value = "2022-05-30 T 21:36:43 UTC"
$ws.Range("B3").Value = 30803.0368329025
$ws.Range("C3").Value = 1.013199
$ws.Range("D3").Value = 1.204574
